$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.699.84'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.510.11'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.83'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.508.63'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.486'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.143'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.56'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.430'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '32.52'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000214'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.104.90'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.508.78'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.642.33'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.117'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.69'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '446.52'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.630'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.653.13'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000127'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.87'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.17%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.49'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.64'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.168'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.45%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.61'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.16'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.86'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.503.70'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.99'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.31'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.37%  '
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '173.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.84%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0894'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.48'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '30.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.882'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.30'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.51'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.63'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.254'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.03%  '
